# Qualis users credentials workbook update
# - Remove the stale "newPT_*" password-reset rows and the old opsadmin
#   password string.
# - Re-point the opsadmin row's password cell to the new "!Qualis1!" value
#   (and drop its mailto hyperlink, since it is no longer an email).
# - Add two new rows: "fund manager" and "invalid user".
# - Refresh the mailto hyperlinks so they keep tracking the right cells,
#   and plain-format everything except the brand-new "invalid user" row
#   (matches the source sheet, where only that row kept the blue/underline
#   Hyperlink look).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Start clean: drop every hyperlink on the sheet; we'll re-create exactly
# the ones that should remain (and add the new one) further down so the
# relationship ids line up with the refreshed content.
$ws.Cells.Hyperlinks.Delete()

# ---- Row 2: superuser (unchanged) ----
$ws.Range("A2").Value = "superuser"
$ws.Range("B2").Value = "autoqualissuperuser@praemium.com"
$ws.Range("C2").Value = "QS@superuserPSS123!"

# ---- Row 3: opsadmin (password replaced, no longer a mailto target) ----
$ws.Range("A3").Value = "opsadmin"
$ws.Range("B3").Value = "autoqualisuser_opsadmin@praemium.com"
$ws.Range("C3").Value = "!Qualis1!"

# ---- Row 4: viewuser (unchanged) ----
$ws.Range("A4").Value = "viewuser"
$ws.Range("B4").Value = "autoqualisuser_view@praemium.com"
$ws.Range("C4").Value = "QS@viewuserPSS123!"

# ---- Row 5: finadviser (unchanged) ----
$ws.Range("A5").Value = "finadviser"
$ws.Range("B5").Value = "autofinadviser@praemium.com"
$ws.Range("C5").Value = "QS@FinAdviser1PSS123!"

# ---- Row 6: forgot password user (new reset-token password) ----
$ws.Range("A6").Value = "forgot password user"
$ws.Range("B6").Value = "qataskdemoaccnt@gmail.com"
$ws.Range("C6").Value = "newPT_724*740"

# ---- Row 7 (new): fund manager ----
$ws.Range("A7").Value = "fund manager"
$ws.Range("B7").Value = "autofundmanager@praemium.com"
$ws.Range("C7").Value = "QS@fundManagerPSS123!"

# ---- Row 8 (new): invalid user ----
$ws.Range("A8").Value = "invalid user"
$ws.Range("B8").Value = "test@test.com"
$ws.Range("C8").Value = "pass12345!."

# Plain-format the whole used range first (clears the legacy custom
# "vertical center" styles off row 6 and the stale Hyperlink look off
# every other linked cell).
$ws.Range("A1:C8").Style = "Normal"

# ---- Rebuild the mailto hyperlinks on the surviving / new email+password cells ----
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:autoqualissuperuser@praemium.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:QS@superuserPSS123!") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:autoqualisuser_opsadmin@praemium.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:autoqualisuser_view@praemium.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:QS@viewuserPSS123!") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:qataskdemoaccnt@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:test@test.com") | Out-Null

# Excel's own Hyperlinks.Add re-applies the blue/underline Hyperlink style;
# strip it back off everywhere except the newly added "invalid user" row,
# which is the only one that keeps the highlighted look in the target sheet.
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Style = "Normal"
$ws.Range("B3").Style = "Normal"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Style = "Normal"
$ws.Range("B6").Style = "Normal"
$ws.Range("B8").Style = "Hyperlink"

# Update the selection to reflect the new last row, matching Excel's
# behaviour of leaving the cursor just past the last entered row.
$ws.Range("A8").Select() | Out-Null
